$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" date field text from 2020/01/26
#    to 2020-05-20 everywhere it appears: the slide master and every slide
#    layout each carry their own Date Placeholder.
# ---------------------------------------------------------------------------
$newDate = "2020-05-20"

function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.HasTextFrame) {
            $isDatePlaceholder = $false
            try {
                if ($shp.PlaceholderFormat.Type -eq 16) {
                    $isDatePlaceholder = $true
                }
            } catch {
                $isDatePlaceholder = $false
            }
            if ($isDatePlaceholder) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master

for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    Update-DatePlaceholder $master.CustomLayouts.Item($j)
}

# ---------------------------------------------------------------------------
# 2) Highlight the "I" callout box (the pin-header annotation on the hero
#    board diagram) to flag that it must be hand assembled: give it a
#    semi-transparent yellow fill instead of no fill.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)
$highlight = $s.Shapes.Item("Rectangle 18")
$highlight.Fill.Visible = $true
$highlight.Fill.Solid()
$highlight.Fill.ForeColor.RGB = 65535
$highlight.Fill.Transparency = 0.75
